# Fix Training Data Issue (#48)
# The BF column holds a "Date" label (BF1) followed by a date string per
# data row (BF2:BF31). The stored value "6-20-2013-14" was off by a day
# because of how the NBA stats site reported the date; the corrected
# representation is the ISO-style "2014-06-20" text value.
#
# NOTE: the date-like text must stay literal text (not get auto-converted
# by Excel into a date serial number), so we briefly force the cells to a
# Text number format before writing the new value, then restore the
# cells' style back to the workbook default ("Normal") to match the
# original (unstyled) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 31
$col = "BF"
$newDate = "2014-06-20"

$rng = $ws.Range($col + $firstDataRow + ":" + $col + $lastDataRow)

# Force text formatting so Excel doesn't reinterpret the string as a date.
$rng.NumberFormat = "@"

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $cell = $ws.Range($col + $row)
    $cell.Value = $newDate
}

# Restore the default "Normal" style so these cells remain unstyled, just
# like they were before the edit.
$rng.Style = "Normal"
